$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 1 values right (A1:E1 -> B1:E1, drop old E1 "12:00 - 12:50"), clear A1
$ws.Range("E1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = $ws.Range("A1").Value()
$ws.Range("A1").ClearContents()

# Replace remaining single-letter slot labels with "SLOT"
$ws.Range("B2").Value = "SLOT"
$ws.Range("I2").Value = "SLOT"
$ws.Range("D3").Value = "SLOT"
$ws.Range("E4").Value = "SLOT"
$ws.Range("C5").Value = "SLOT"
$ws.Range("H5").Value = "SLOT"
$ws.Range("I6").Value = "SLOT"

Write-Host "done"
